$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.429192543029785
$ws.Range("B1").Value = 1.597786068916321
$ws.Range("C1").Value = 6.532034397125244
$ws.Range("D1").Value = 1.635138988494873
$ws.Range("E1").Value = 0.9716724753379822
